$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 4930
$ws.Range("I32").Value = 4471.4287
$ws.Range("J32").Value = 6000
$ws.Range("K32").Value = 4471.4287
$ws.Range("L32").Value = 6000
$ws.Range("M32").Value = -4145.4287
$ws.Range("N32").Value = -6652

$ws.Range("H112").Value = 1660.4736
$ws.Range("I112").Value = 1299.8
$ws.Range("J112").Value = 1789.2858
$ws.Range("K112").Value = 3899.4
$ws.Range("L112").Value = 5367.857400000001
$ws.Range("M112").Value = -2791.4
$ws.Range("N112").Value = -7583.857400000001

$ws.Range("H132").Value = 2969.1228
$ws.Range("I132").Value = 2772.6274
$ws.Range("J132").Value = 4639.3335
$ws.Range("K132").Value = 8317.8822
$ws.Range("L132").Value = 13918.0005
$ws.Range("M132").Value = -5787.8822
$ws.Range("N132").Value = -18978.0005

$ws.Range("H138").Value = 2844.4443
$ws.Range("I138").Value = 1921.9
$ws.Range("J138").Value = 3387.1177
$ws.Range("K138").Value = 5765.700000000001
$ws.Range("L138").Value = 10161.3531
$ws.Range("M138").Value = -625.7000000000007
$ws.Range("N138").Value = -20441.3531

$ws.Range("H141").Value = 4962.2
$ws.Range("I141").Value = 4764.3076
$ws.Range("J141").Value = 6248.5
$ws.Range("K141").Value = 14292.9228
$ws.Range("L141").Value = 18745.5
$ws.Range("M141").Value = -9112.9228
$ws.Range("N141").Value = -29105.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 6077.909
$ws.Range("I31").Value = 6077.909
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 6077.909
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -5783.909

$ws.Range("H34").Value = 179750
$ws.Range("I34").Value = 59500
$ws.Range("J34").Value = 300000
$ws.Range("K34").Value = 59500
$ws.Range("L34").Value = 300000
$ws.Range("M34").Value = -59229
$ws.Range("N34").Value = -300542

$ws.Range("H102").Value = 12792.972
$ws.Range("I102").Value = 17269.05
$ws.Range("J102").Value = 6824.8667
$ws.Range("K102").Value = 17269.05
$ws.Range("L102").Value = 6824.8667
$ws.Range("M102").Value = -15647.05
$ws.Range("N102").Value = -10068.8667

$ws.Range("H122").Value = 628454.9
$ws.Range("I122").Value = 3161.697
$ws.Range("J122").Value = 2004099.9
$ws.Range("K122").Value = 9485.091
$ws.Range("L122").Value = 6012299.699999999
$ws.Range("M122").Value = -7035.091
$ws.Range("N122").Value = -6017199.699999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2135.9033
$ws.Range("I94").Value = 1167.1
$ws.Range("J94").Value = 3897.3635
$ws.Range("K94").Value = 1167.1
$ws.Range("L94").Value = 3897.3635
$ws.Range("M94").Value = -716.0999999999999
$ws.Range("N94").Value = -4799.363499999999

$ws.Range("H102").Value = 24999.5
$ws.Range("I102").Value = 14999.5
$ws.Range("J102").Value = 34999.5
$ws.Range("K102").Value = 14999.5
$ws.Range("L102").Value = 34999.5
$ws.Range("M102").Value = -11754.5
$ws.Range("N102").Value = -41489.5

$ws.Range("H107").Value = 2025
$ws.Range("I107").Value = 2080.2632
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 2080.2632
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = -160.2631999999999
$ws.Range("N107").Value = -5340

$ws.Range("H134").Value = 15790.363
$ws.Range("I134").Value = 20671.428
$ws.Range("J134").Value = 7248.5
$ws.Range("K134").Value = 62014.284
$ws.Range("L134").Value = 21745.5
$ws.Range("M134").Value = -59479.284
$ws.Range("N134").Value = -26815.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4813.2046
$ws.Range("I31").Value = 4879.514
$ws.Range("J31").Value = 4555.3335
$ws.Range("K31").Value = 4879.514
$ws.Range("L31").Value = 4555.3335
$ws.Range("M31").Value = -4584.514
$ws.Range("N31").Value = -5145.3335

$ws.Range("H34").Value = 4813.2046
$ws.Range("I34").Value = 4879.514
$ws.Range("J34").Value = 4555.3335
$ws.Range("K34").Value = 4879.514
$ws.Range("L34").Value = 4555.3335
$ws.Range("M34").Value = -4677.514
$ws.Range("N34").Value = -4959.3335

$ws.Range("H52").Value = 69500
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 69500
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 69500
$ws.Range("N52").Value = -70088

$ws.Range("H58").Value = 3515.6924
$ws.Range("I58").Value = 4116.077
$ws.Range("J58").Value = 2915.3076
$ws.Range("K58").Value = 4116.077
$ws.Range("L58").Value = 2915.3076
$ws.Range("M58").Value = -3913.077
$ws.Range("N58").Value = -3321.3076

$ws.Range("H99").Value = 316888.3
$ws.Range("I99").Value = 627589.25
$ws.Range("J99").Value = 6187.375
$ws.Range("K99").Value = 627589.25
$ws.Range("L99").Value = 6187.375
$ws.Range("M99").Value = -626091.25
$ws.Range("N99").Value = -9183.375

$ws.Range("H122").Value = 16411.375
$ws.Range("I122").Value = 31073
$ws.Range("J122").Value = 1749.75
$ws.Range("K122").Value = 93219
$ws.Range("L122").Value = 5249.25
$ws.Range("M122").Value = -90769
$ws.Range("N122").Value = -10149.25

$ws.Range("H126").Value = 316888.3
$ws.Range("I126").Value = 627589.25
$ws.Range("J126").Value = 6187.375
$ws.Range("K126").Value = 1882767.75
$ws.Range("L126").Value = 18562.125
$ws.Range("M126").Value = -1880297.75
$ws.Range("N126").Value = -23502.125

$ws.Range("H130").Value = 110000
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 110000
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 110000
$ws.Range("N130").Value = -120040

$ws.Range("H136").Value = 3515.6924
$ws.Range("I136").Value = 4116.077
$ws.Range("J136").Value = 2915.3076
$ws.Range("K136").Value = 12348.231
$ws.Range("L136").Value = 8745.9228
$ws.Range("M136").Value = -9798.231
$ws.Range("N136").Value = -13845.9228

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 10528953
$ws.Range("I131").Value = 83340390
$ws.Range("J131").Value = 1998.8313
$ws.Range("K131").Value = 250021170
$ws.Range("L131").Value = 5996.4939
$ws.Range("M131").Value = -250016130
$ws.Range("N131").Value = -16076.4939

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").ClearContents()
$ws.Range("N68").Value = 0

$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").ClearContents()
$ws.Range("N71").Value = 0

$ws.Range("H74").Value = 50119
$ws.Range("I74").Value = 50119
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 50119
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -49183

$ws.Range("H77").Value = 50119
$ws.Range("I77").Value = 50119
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 150357
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -145677

$ws.Range("H80").Value = 11742.789
$ws.Range("I80").Value = 15262.75
$ws.Range("J80").Value = 5708.5713
$ws.Range("K80").Value = 15262.75
$ws.Range("L80").Value = 5708.5713
$ws.Range("M80").Value = -14264.75
$ws.Range("N80").Value = -7704.5713

$ws.Range("H83").Value = 11742.789
$ws.Range("I83").Value = 15262.75
$ws.Range("J83").Value = 5708.5713
$ws.Range("K83").Value = 76313.75
$ws.Range("L83").Value = 28542.8565
$ws.Range("M83").Value = -71321.75
$ws.Range("N83").Value = -38526.85649999999

$ws.Range("H103").Value = 35000
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 35000
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 35000
$ws.Range("N103").Value = -37344

$ws.Range("H136").Value = 15298.588
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 15298.588
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 45895.764
$ws.Range("N136").Value = -50995.764

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 24176.363
$ws.Range("I7").Value = 33420.285
$ws.Range("J7").Value = 7999.5
$ws.Range("K7").Value = 33420.285
$ws.Range("L7").Value = 7999.5
$ws.Range("M7").Value = -33308.285
$ws.Range("N7").Value = -8223.5

$ws.Range("H40").Value = 29283.37
$ws.Range("I40").Value = 38115.832
$ws.Range("J40").Value = 14142
$ws.Range("K40").Value = 38115.832
$ws.Range("L40").Value = 14142
$ws.Range("M40").Value = -37979.832
$ws.Range("N40").Value = -14414

$ws.Range("H61").Value = 4645
$ws.Range("I61").Value = 1141.909
$ws.Range("J61").Value = 11067.333
$ws.Range("K61").Value = 1141.909
$ws.Range("L61").Value = 11067.333
$ws.Range("M61").Value = -939.9090000000001
$ws.Range("N61").Value = -11471.333

$ws.Range("H113").Value = 4645
$ws.Range("I113").Value = 1141.909
$ws.Range("J113").Value = 11067.333
$ws.Range("K113").Value = 1141.909
$ws.Range("L113").Value = 11067.333
$ws.Range("M113").Value = 1028.091
$ws.Range("N113").Value = -15407.333

$ws.Range("H126").Value = 24176.363
$ws.Range("I126").Value = 33420.285
$ws.Range("J126").Value = 7999.5
$ws.Range("K126").Value = 100260.855
$ws.Range("L126").Value = 23998.5
$ws.Range("M126").Value = -97790.85500000001
$ws.Range("N126").Value = -28938.5

$ws.Range("H136").Value = 7934.3335
$ws.Range("I136").Value = 13999
$ws.Range("J136").Value = 7383
$ws.Range("K136").Value = 41997
$ws.Range("L136").Value = 22149
$ws.Range("M136").Value = -39447
$ws.Range("N136").Value = -27249

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 38211.57
$ws.Range("I75").Value = 27500
$ws.Range("J75").Value = 39996.832
$ws.Range("K75").Value = 27500
$ws.Range("L75").Value = 39996.832
$ws.Range("M75").Value = -26564
$ws.Range("N75").Value = -41868.832

$ws.Range("H78").Value = 38211.57
$ws.Range("I78").Value = 27500
$ws.Range("J78").Value = 39996.832
$ws.Range("K78").Value = 82500
$ws.Range("L78").Value = 119990.496
$ws.Range("M78").Value = -77820
$ws.Range("N78").Value = -129350.496

$ws.Range("H80").Value = 29997
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 29997
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 29997
$ws.Range("N80").Value = -31993

$ws.Range("H83").Value = 29997
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 29997
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 89991
$ws.Range("N83").Value = -99975

$ws.Range("H132").Value = 15960
$ws.Range("I132").Value = 26452.846
$ws.Range("J132").Value = 6866.2
$ws.Range("K132").Value = 79358.538
$ws.Range("L132").Value = 20598.6
$ws.Range("M132").Value = -76828.538
$ws.Range("N132").Value = -25658.6

$ws.Range("H136").Value = 411681.56
$ws.Range("I136").Value = 575134.8
$ws.Range("J136").Value = 10478.091
$ws.Range("K136").Value = 1725404.4
$ws.Range("L136").Value = 31434.273
$ws.Range("M136").Value = -1722854.4
$ws.Range("N136").Value = -36534.273
